# Generate Report for Handoff
#
# The "e0390748-aa6c-4338-96c4-c0d594e4314f" entry and the
# "8e2b4336-5735-4c5d-bffc-e21e7f60d885" entry swap places (e0390748 moves to
# row 2, 8e2b4336 moves to row 3) on the Overview sheet and on each
# per-language detail sheet (zh-cn, de-de). In addition, the 8e2b4336 entry
# (now in row 3) is marked as ready for a new handoff again, with fresh
# handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("A2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.md"
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"
$overview.Range("D2").Value = "2016-30-14 04:30:10"

$overview.Range("A3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.md"
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-30-14 04:30:59"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("A2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.md"
$zhcn.Range("B2").Value = ".md"
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("D2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.zh-cn.xlf"
$zhcn.Range("E2").Value = "2016-03-14 04:30:06"
$zhcn.Range("F2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.md"
$zhcn.Range("G2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-03-14 04:30:29"
$zhcn.Range("I2").Value = "Include"

$zhcn.Range("A3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.md"
$zhcn.Range("B3").Value = ".md"
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-14 04:30:56"
$zhcn.Range("F3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.md"
$zhcn.Range("G3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-03-14 04:30:29"
$zhcn.Range("I3").Value = "Include"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("A2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.md"
$dede.Range("B2").Value = ".md"
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("D2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.de-de.xlf"
$dede.Range("E2").Value = "2016-03-14 04:30:10"
$dede.Range("F2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.md"
$dede.Range("G2").Value = "e0390748-aa6c-4338-96c4-c0d594e4314f.6f38887ef8fb072b659eafcafeb8544e61d5db31.de-de.xlf"
$dede.Range("H2").Value = "2016-03-14 04:30:34"
$dede.Range("I2").Value = "Include"

$dede.Range("A3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.md"
$dede.Range("B3").Value = ".md"
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.de-de.xlf"
$dede.Range("E3").Value = "2016-03-14 04:30:59"
$dede.Range("F3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.md"
$dede.Range("G3").Value = "8e2b4336-5735-4c5d-bffc-e21e7f60d885.2a5251d41fd3c63e7892720e50d02aaef985fb9a.de-de.xlf"
$dede.Range("H3").Value = "2016-03-14 04:30:34"
$dede.Range("I3").Value = "Include"
